$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.993.29'
$ws.Range('E2').Value = '  -6.05%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.431.99'
$ws.Range('E3').Value = '  -7.76%  '

$ws.Range('E4').Value = '  +0.27%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '551.15'
$ws.Range('E5').Value = '  -10.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.25'
$ws.Range('E6').Value = '  -5.66%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.432.72'
$ws.Range('E7').Value = '  -7.72%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  -6.97%  '

$ws.Range('E9').Value = '  +0.27%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.635'
$ws.Range('E10').Value = '  -12.26%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.138'
$ws.Range('E11').Value = '  -14.07%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '50.59'
$ws.Range('E12').Value = '  -16.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('E13').Value = '  -15.17%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.30'
$ws.Range('E14').Value = '  -12.48%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.986.02'
$ws.Range('E15').Value = '  -7.76%  '

$ws.Range('E16').Value = '  -1.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.450.51'
$ws.Range('E17').Value = '  -7.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.778.26'
$ws.Range('E18').Value = '  -6.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.53'
$ws.Range('E19').Value = '  -9.68%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.52'
$ws.Range('E20').Value = '  -10.76%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.01'
$ws.Range('E21').Value = '  -11.43%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '373.47'
$ws.Range('E22').Value = '  -9.35%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.04'
$ws.Range('E23').Value = '  -11.86%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.02'
$ws.Range('E24').Value = '  -8.29%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.49'
$ws.Range('E25').Value = '  -3.34%  '

$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.98'
$ws.Range('E26').Value = '  -1.13%  '

$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.75'
$ws.Range('E27').Value = '  -9.78%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.69'
$ws.Range('E28').Value = '  -8.94%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.35'
$ws.Range('E29').Value = '  -11.94%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.38'
$ws.Range('E30').Value = '  -13.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '29.70'
$ws.Range('E31').Value = '  -10.20%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.86'
$ws.Range('E32').Value = '  -10.31%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '599.89'
$ws.Range('E33').Value = '  -4.89%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.61'
$ws.Range('E34').Value = '  -8.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '62.35'
$ws.Range('E35').Value = '  -5.11%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.108'
$ws.Range('E36').Value = '  -12.96%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.54'
$ws.Range('E37').Value = '  -13.62%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.383'
$ws.Range('E39').Value = '  -7.54%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0702'
$ws.Range('E41').Value = '  -14.93%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.127'
$ws.Range('E42').Value = '  -10.15%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.897.70'
$ws.Range('E43').Value = '  +0.52%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.67'
$ws.Range('E44').Value = '  -12.52%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.12'
$ws.Range('E45').Value = '  +0.32%  '

$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.38'
$ws.Range('E46').Value = '  -9.26%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0387'
$ws.Range('E47').Value = '  -13.25%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.125'
$ws.Range('E48').Value = '  -10.33%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '136.78'
$ws.Range('E49').Value = '  -3.74%  '

$ws.Range('E50').Value = '  -10.83%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.07'
$ws.Range('E51').Value = '  -12.13%  '
